$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 425
$firstRow = 2

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45182
